$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "After H6 SpiderCrawl..." row), shifting all rows below up by one.
$ws.Rows.Item(2).Delete()

# Add the new note to what is now row 2 ("locked doors, better doors?, door handles anim")
$ws.Range("E2").Value = '"Nothing in there", "Locked?"'
